# Auto update: 2025-11-29 18:42:09
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Riot Platforms, Inc. (RIOT) ---
$ws.Range("E2").Value = 44.9
$ws.Range("G2").Value = 50
$ws.Range("K2").Value = 60.8
$ws.Range("M2").Value = "📈 매수 관찰 구간입니다."
$ws.Range("N2").Value = 85.87127175646313
$ws.Range("O2").Value = "🟢 완화적 (상승 우위)"

# --- Row 3: now Bitcoin USD (BTC-USD) with refreshed data ---
$ws.Range("B3").Value = "Bitcoin USD"
$ws.Range("C3").Value = "BTC-USD"
$ws.Range("D3").Value = 90939.73
$ws.Range("E3").Value = 42.4
$ws.Range("F3").Value = 4.76
$ws.Range("G3").Value = 40
$ws.Range("H3").Value = 56
$ws.Range("I3").Value = 43
$ws.Range("J3").Value = 43
$ws.Range("K3").Value = 55
$ws.Range("M3").Value = "⛔ 관망하십시오."
$ws.Range("N3").Value = 85.87127175646313
$ws.Range("O3").Value = "🟢 완화적 (상승 우위)"

# --- Row 4: now Coinbase Global, Inc. (COIN) ---
$ws.Range("B4").Value = "Coinbase Global, Inc."
$ws.Range("C4").Value = "COIN"
$ws.Range("D4").Value = 272.82
$ws.Range("E4").Value = 35.9
$ws.Range("F4").Value = 14.55
$ws.Range("G4").Value = 20
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 53
$ws.Range("J4").Value = 46
$ws.Range("K4").Value = 53
$ws.Range("M4").Value = "⛔ 관망하십시오."
$ws.Range("N4").Value = 85.87127175646313
$ws.Range("O4").Value = "🟢 완화적 (상승 우위)"

# --- Row 5: MARA Holdings, Inc. (MARA) ---
$ws.Range("E5").Value = 26.1
$ws.Range("M5").Value = "⛔ 관망하십시오."
$ws.Range("N5").Value = 85.87127175646313
$ws.Range("O5").Value = "🟢 완화적 (상승 우위)"

# --- Row 6: Strategy Inc (MSTR) ---
$ws.Range("E6").Value = 21.7
$ws.Range("M6").Value = "⛔ 관망하십시오."
$ws.Range("N6").Value = 85.87127175646313
$ws.Range("O6").Value = "🟢 완화적 (상승 우위)"
